$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "Цены"
$ws2 = $wb.Worksheets.Item(2)   # "Продажи 2"

# ---------------------------------------------------------------------------
# 1) New sheet-scoped ("local") defined names "Курс" / "Цена" on "Продажи 2",
#    alongside the pre-existing workbook-scoped ones.
#    (Names.Add can't take a non-ASCII name directly in this runtime, so we
#    add with a placeholder name, then rename.)
# ---------------------------------------------------------------------------
$localTsena = $ws2.Names.Add('TmpTsenaLocal', '=''Продажи 2''!$G$1:$G$16')
$localTsena.Name = 'Цена'

$localKurs = $ws2.Names.Add('TmpKursLocal', '=''Продажи 2''!$H$1:$H$16')
$localKurs.Name = 'Курс'

# ---------------------------------------------------------------------------
# 2) Sheet1 ("Цены"): F2:F16 get centre-aligned; F7 text "шт." -> "метр"
# ---------------------------------------------------------------------------
$ws1.Range("F2:F16").HorizontalAlignment = -4108
$ws1.Range("F7").Value = "метр"

# ---------------------------------------------------------------------------
# 3) Sheet1 & Sheet2: G14 empty/5 -> 13 (copy number format from G12 so the
#    style matches what Excel would reuse for a freshly-entered number)
# ---------------------------------------------------------------------------
$ws1.Range("G14").Value = 13
$ws1.Range("G12").Copy()
$ws1.Range("G14").PasteSpecial(-4122)

$ws2.Range("G14").Value = 13
$ws2.Range("G12").Copy()
$ws2.Range("G14").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 4) Sheet2 ("Продажи 2"): I2:I16 formulas switch from H*G to the named-range
#    based formula, each one an explicit (non-shared) formula.
# ---------------------------------------------------------------------------
$ws2.Range("I2").Formula  = "='Продажи 2'!Курс*'Продажи 2'!Цена"
$ws2.Range("I3").Formula  = "='Продажи 2'!Курс*'Продажи 2'!Цена"
$ws2.Range("I4").Formula  = "='Продажи 2'!Курс*'Продажи 2'!Цена"
$ws2.Range("I5").Formula  = "='Продажи 2'!Курс*'Продажи 2'!Цена"
$ws2.Range("I6").Formula  = "='Продажи 2'!Курс*'Продажи 2'!Цена"
$ws2.Range("I7").Formula  = "='Продажи 2'!Курс*'Продажи 2'!Цена"
$ws2.Range("I8").Formula  = "='Продажи 2'!Курс*'Продажи 2'!Цена"
$ws2.Range("I9").Formula  = "='Продажи 2'!Курс*'Продажи 2'!Цена"
$ws2.Range("I10").Formula = "='Продажи 2'!Курс*'Продажи 2'!Цена"
$ws2.Range("I11").Formula = "='Продажи 2'!Курс*'Продажи 2'!Цена"
$ws2.Range("I12").Formula = "='Продажи 2'!Курс*'Продажи 2'!Цена"
$ws2.Range("I13").Formula = "='Продажи 2'!Курс*'Продажи 2'!Цена"
$ws2.Range("I14").Formula = "='Продажи 2'!Курс*'Продажи 2'!Цена"
$ws2.Range("I15").Formula = "='Продажи 2'!Курс*'Продажи 2'!Цена"
$ws2.Range("I16").Formula = "='Продажи 2'!Курс*'Продажи 2'!Цена"

# ---------------------------------------------------------------------------
# 5) Sheet1 page setup: drop verticalDpi (touching PageSetup rewrites the
#    element without it).
# ---------------------------------------------------------------------------
$ws1.PageSetup.Orientation = 1

# ---------------------------------------------------------------------------
# 6) Selections / view state
# ---------------------------------------------------------------------------
$ws1.Range("G2").Select()
$ws2.Range("I16").Select()
$ws2.Activate()
$excel.ActiveWindow.DisplayFormulas = $true
